$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - column F ("想去人数" / want-to-go count) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 636
$ws1.Range("F3").Value = 688
$ws1.Range("F4").Value = 941
$ws1.Range("F5").Value = 711
$ws1.Range("F7").Value = 394
$ws1.Range("F9").Value = 126
$ws1.Range("F10").Value = 1199
$ws1.Range("F12").Value = 380
$ws1.Range("F13").Value = 503
$ws1.Range("F15").Value = 8
$ws1.Range("F16").Value = 426
$ws1.Range("F17").Value = 345
$ws1.Range("F21").Value = 70
$ws1.Range("F22").Value = 567
$ws1.Range("F23").Value = 26
$ws1.Range("F24").Value = 725

# Sheet "演出" (Performances) - column F updates
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 84
$ws2.Range("F4").Value = 314
$ws2.Range("F6").Value = 19
$ws2.Range("F9").Value = 219
$ws2.Range("F10").Value = 48
$ws2.Range("F11").Value = 24
$ws2.Range("F13").Value = 86

# Sheet "全部类型" (All types) - column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 84
$ws4.Range("F4").Value = 636
$ws4.Range("F6").Value = 314
$ws4.Range("F7").Value = 688
$ws4.Range("F8").Value = 941
$ws4.Range("F9").Value = 711
$ws4.Range("F11").Value = 394
$ws4.Range("F13").Value = 126
$ws4.Range("F14").Value = 1199
$ws4.Range("F17").Value = 19
$ws4.Range("F18").Value = 380
$ws4.Range("F19").Value = 503
$ws4.Range("F22").Value = 8
$ws4.Range("F23").Value = 426
$ws4.Range("F25").Value = 345
$ws4.Range("F28").Value = 219
$ws4.Range("F29").Value = 48
$ws4.Range("F31").Value = 24
$ws4.Range("F33").Value = 86
$ws4.Range("F34").Value = 70
$ws4.Range("F35").Value = 567
$ws4.Range("F36").Value = 26
$ws4.Range("F37").Value = 725
